# Weekly update: insert the newest Maracuyá price record at row 83
# (Vega Modelo de Temuco / Región de Arica y Parinacota), pushing the
# existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 83, shifting rows 83:112 down to 84:113.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly record.
$ws.Cells.Item(83, 1).Value2 = 10
$ws.Cells.Item(83, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value2 = "La Araucanía"
$ws.Cells.Item(83, 4).Value2 = 45146
$ws.Cells.Item(83, 5).Value2 = 9
$ws.Cells.Item(83, 6).Value2 = "Fruta"
$ws.Cells.Item(83, 7).Value2 = 100108
$ws.Cells.Item(83, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(83, 9).Value2 = 100108003
$ws.Cells.Item(83, 10).Value2 = "Maracuyá"
$ws.Cells.Item(83, 11).Value2 = "Sin especificar"
$ws.Cells.Item(83, 12).Value2 = "Primera"
$ws.Cells.Item(83, 13).Value2 = 30
$ws.Cells.Item(83, 14).Value2 = 42000
$ws.Cells.Item(83, 15).Value2 = 42000
$ws.Cells.Item(83, 16).Value2 = 42000
$ws.Cells.Item(83, 17).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(83, 18).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(83, 19).Value2 = 2333
$ws.Cells.Item(83, 20).Value2 = 18
